# Update LR-pair TPM-based metrics (ligand/receptor expression values,
# specificities, and edge weights) with newly computed TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.727484333333333
$ws.Range("H2").Value = 8.182453000000001
$ws.Range("I2").Value = 0.03096049453772388
$ws.Range("J2").Value = 0.03096049453772388
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 166.4900249556675
$ws.Range("R2").Value = 1498.410224601007
$ws.Range("S2").Value = 0.006327128338333687
$ws.Range("T2").Value = 0.006327128338333687
$ws.Range("G3").Value = 2.727484333333333
$ws.Range("H3").Value = 8.182453000000001
$ws.Range("I3").Value = 0.03096049453772388
$ws.Range("J3").Value = 0.03096049453772388
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 289.9710404216993
$ws.Range("R3").Value = 2609.739363795294
$ws.Range("S3").Value = 0.01101978324309082
$ws.Range("T3").Value = 0.01101978324309082
$ws.Range("G4").Value = 2.727484333333333
$ws.Range("H4").Value = 8.182453000000001
$ws.Range("I4").Value = 0.03096049453772388
$ws.Range("J4").Value = 0.03096049453772388
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 358.2234538215871
$ws.Range("R4").Value = 3224.011084394284
$ws.Range("S4").Value = 0.01361358295629937
$ws.Range("T4").Value = 0.01361358295629938
$ws.Range("I5").Value = 0.5986009007423507
$ws.Range("J5").Value = 0.5986009007423507
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 3218.97567823559
$ws.Range("R5").Value = 28970.78110412031
$ws.Range("S5").Value = 0.1223308858269109
$ws.Range("T5").Value = 0.1223308858269109
$ws.Range("I6").Value = 0.5986009007423507
$ws.Range("J6").Value = 0.5986009007423507
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("Q6").Value = 5606.400303914094
$ws.Range("R6").Value = 50457.60273522685
$ws.Range("S6").Value = 0.2130602974465464
$ws.Range("T6").Value = 0.2130602974465465
$ws.Range("I7").Value = 0.5986009007423507
$ws.Range("J7").Value = 0.5986009007423507
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 6926.016051305696
$ws.Range("R7").Value = 62334.14446175127
$ws.Range("S7").Value = 0.2632097174688934
$ws.Range("T7").Value = 0.2632097174688934
$ws.Range("G8").Value = 32.63402300000001
$ws.Range("H8").Value = 97.90206900000001
$ws.Range("I8").Value = 0.3704386047199253
$ws.Range("J8").Value = 0.3704386047199253
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 1992.033185038946
$ws.Range("R8").Value = 17928.29866535052
$ws.Range("S8").Value = 0.07570333189220885
$ws.Range("T8").Value = 0.07570333189220885
$ws.Range("G9").Value = 32.63402300000001
$ws.Range("H9").Value = 97.90206900000001
$ws.Range("I9").Value = 0.3704386047199253
$ws.Range("J9").Value = 0.3704386047199253
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 3469.468728676718
$ws.Range("R9").Value = 31225.21855809047
$ws.Range("S9").Value = 0.131850385138799
$ws.Range("T9").Value = 0.131850385138799
$ws.Range("G10").Value = 32.63402300000001
$ws.Range("H10").Value = 97.90206900000001
$ws.Range("I10").Value = 0.3704386047199253
$ws.Range("J10").Value = 0.3704386047199253
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 4286.100670967415
$ws.Range("R10").Value = 38574.90603870674
$ws.Range("S10").Value = 0.1628848876889174
$ws.Range("T10").Value = 0.1628848876889175
